# Auto-generated edit script applying numeric corrections to the leve-profit
# tables on each job sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Source data refreshed by the scheduled Sheets runner; this script replays
# the resulting cell-level deltas via the Excel COM object model.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 83427.836
$ws.Range("J6").Value = 50
$ws.Range("L6").Value = 150
$ws.Range("N6").Value = -374
$ws.Range("H9").Value = 98.08108
$ws.Range("J9").Value = 70.5
$ws.Range("L9").Value = 70.5
$ws.Range("N9").Value = -408.5
$ws.Range("H28").Value = 2696.6667
$ws.Range("I28").Value = 2676.6667
$ws.Range("J28").Value = 2756.6667
$ws.Range("K28").Value = 2676.6667
$ws.Range("L28").Value = 2756.6667
$ws.Range("M28").Value = -2191.6667
$ws.Range("N28").Value = -3726.6667
$ws.Range("H38").Value = 5547.636
$ws.Range("I38").Value = 204.8
$ws.Range("J38").Value = 10000
$ws.Range("K38").Value = 614.4000000000001
$ws.Range("L38").Value = 30000
$ws.Range("M38").Value = -242.4000000000001
$ws.Range("N38").Value = -30744
$ws.Range("H76").Value = 17671.4
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 17671.4
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 17671.4
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -18301.4
$ws.Range("H79").Value = 17671.4
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 17671.4
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 17671.4
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -19855.4
$ws.Range("H88").Value = 2851.1428
$ws.Range("I88").Value = 2700
$ws.Range("J88").Value = 3229
$ws.Range("K88").Value = 2700
$ws.Range("L88").Value = 3229
$ws.Range("M88").Value = -2294
$ws.Range("N88").Value = -4041
$ws.Range("H91").Value = 2851.1428
$ws.Range("I91").Value = 2700
$ws.Range("J91").Value = 3229
$ws.Range("K91").Value = 2700
$ws.Range("L91").Value = 3229
$ws.Range("M91").Value = -1296
$ws.Range("N91").Value = -6037
$ws.Range("H112").Value = 3823.7778
$ws.Range("J112").Value = 3995.3333
$ws.Range("L112").Value = 11985.9999
$ws.Range("N112").Value = -14201.9999
$ws.Range("H116").Value = 4365.6665
$ws.Range("J116").Value = 5500
$ws.Range("L116").Value = 5500
$ws.Range("N116").Value = -12384
$ws.Range("H118").Value = 454.1111
$ws.Range("I118").Value = 454.1111
$ws.Range("K118").Value = 1362.3333
$ws.Range("M118").Value = 294.6667
$ws.Range("H132").Value = 4958.0713
$ws.Range("I132").Value = 5068.6484
$ws.Range("K132").Value = 15205.9452
$ws.Range("M132").Value = -12675.9452
$ws.Range("H138").Value = 2534.3965
$ws.Range("I138").Value = 1331.2188
$ws.Range("J138").Value = 4015.2307
$ws.Range("K138").Value = 3993.6564
$ws.Range("L138").Value = 12045.6921
$ws.Range("M138").Value = 1146.3436
$ws.Range("N138").Value = -22325.6921

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9909.978999999999
$ws.Range("I61").Value = 6907.6
$ws.Range("K61").Value = 6907.6
$ws.Range("M61").Value = -6695.6
$ws.Range("H74").Value = 3790.95
$ws.Range("I74").Value = 3340.1
$ws.Range("J74").Value = 4692.65
$ws.Range("K74").Value = 3340.1
$ws.Range("L74").Value = 4692.65
$ws.Range("M74").Value = -2466.1
$ws.Range("N74").Value = -6440.65
$ws.Range("H77").Value = 3790.95
$ws.Range("I77").Value = 3340.1
$ws.Range("J77").Value = 4692.65
$ws.Range("K77").Value = 16700.5
$ws.Range("L77").Value = 23463.25
$ws.Range("M77").Value = -12332.5
$ws.Range("N77").Value = -32199.25
$ws.Range("H122").Value = 1292.2778
$ws.Range("J122").Value = 1853.6666
$ws.Range("L122").Value = 5560.9998
$ws.Range("N122").Value = -10460.9998
$ws.Range("H132").Value = 3506.102
$ws.Range("I132").Value = 3628.6
$ws.Range("K132").Value = 10885.8
$ws.Range("M132").Value = -8355.799999999999
$ws.Range("H136").Value = 9909.978999999999
$ws.Range("I136").Value = 6907.6
$ws.Range("K136").Value = 20722.8
$ws.Range("M136").Value = -18172.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 59769.5
$ws.Range("I2").Value = 50000
$ws.Range("J2").Value = 69539
$ws.Range("K2").Value = 50000
$ws.Range("L2").Value = 69539
$ws.Range("M2").Value = -49887
$ws.Range("N2").Value = -69765
$ws.Range("H141").Value = 199999
$ws.Range("J141").Value = 199999
$ws.Range("L141").Value = 199999
$ws.Range("N141").Value = -210359

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3675.65
$ws.Range("I16").Value = 3348.5454
$ws.Range("K16").Value = 3348.5454
$ws.Range("M16").Value = -3061.5454
$ws.Range("H31").Value = 6219.1055
$ws.Range("I31").Value = 4044.5
$ws.Range("J31").Value = 7222.769
$ws.Range("K31").Value = 4044.5
$ws.Range("L31").Value = 7222.769
$ws.Range("M31").Value = -3749.5
$ws.Range("N31").Value = -7812.769
$ws.Range("H34").Value = 6219.1055
$ws.Range("I34").Value = 4044.5
$ws.Range("J34").Value = 7222.769
$ws.Range("K34").Value = 4044.5
$ws.Range("L34").Value = 7222.769
$ws.Range("M34").Value = -3842.5
$ws.Range("N34").Value = -7626.769
$ws.Range("H74").Value = 39405.145
$ws.Range("J74").Value = 39405.145
$ws.Range("L74").Value = 39405.145
$ws.Range("N74").Value = -41153.145
$ws.Range("H77").Value = 39405.145
$ws.Range("J77").Value = 39405.145
$ws.Range("L77").Value = 118215.435
$ws.Range("N77").Value = -126951.435
$ws.Range("H110").Value = 79600
$ws.Range("J110").Value = 79600
$ws.Range("L110").Value = 79600
$ws.Range("N110").Value = -87780
$ws.Range("H113").Value = 3675.65
$ws.Range("I113").Value = 3348.5454
$ws.Range("K113").Value = 3348.5454
$ws.Range("M113").Value = -1178.5454
$ws.Range("H132").Value = 4268.5454
$ws.Range("I132").Value = 4377.4736
$ws.Range("K132").Value = 13132.4208
$ws.Range("M132").Value = -10602.4208
$ws.Range("H134").Value = 4613.2812
$ws.Range("I134").Value = 3553.08
$ws.Range("K134").Value = 10659.24
$ws.Range("M134").Value = -8124.24

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 547.5
$ws.Range("J34").Value = 675
$ws.Range("L34").Value = 2025
$ws.Range("N34").Value = -2193
$ws.Range("H39").Value = 1533.3334
$ws.Range("J39").Value = 1533.3334
$ws.Range("L39").Value = 4600.0002
$ws.Range("N39").Value = -5188.0002
$ws.Range("H55").Value = 914.8889
$ws.Range("J55").Value = 980
$ws.Range("L55").Value = 2940
$ws.Range("N55").Value = -3294

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5369.143
$ws.Range("I80").Value = 3731.25
$ws.Range("J80").Value = 7553
$ws.Range("K80").Value = 3731.25
$ws.Range("L80").Value = 7553
$ws.Range("M80").Value = -2733.25
$ws.Range("N80").Value = -9549
$ws.Range("H83").Value = 5369.143
$ws.Range("I83").Value = 3731.25
$ws.Range("J83").Value = 7553
$ws.Range("K83").Value = 18656.25
$ws.Range("L83").Value = 37765
$ws.Range("M83").Value = -13664.25
$ws.Range("N83").Value = -47749
$ws.Range("H97").Value = 4990.4
$ws.Range("I97").Value = 1702.625
$ws.Range("J97").Value = 8747.857
$ws.Range("K97").Value = 1702.625
$ws.Range("L97").Value = 8747.857
$ws.Range("M97").Value = -1206.625
$ws.Range("N97").Value = -9739.857
$ws.Range("H99").Value = 9460.875
$ws.Range("I99").Value = 8669.571
$ws.Range("K99").Value = 8669.571
$ws.Range("M99").Value = -6423.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5772.091
$ws.Range("I46").Value = 5772.091
$ws.Range("K46").Value = 5772.091
$ws.Range("M46").Value = -5584.091
$ws.Range("H101").Value = 100000
$ws.Range("J101").Value = 100000
$ws.Range("L101").Value = 100000
$ws.Range("N101").Value = -106490
$ws.Range("H132").Value = 8466.919
$ws.Range("I132").Value = 7216.952
$ws.Range("K132").Value = 21650.856
$ws.Range("M132").Value = -19120.856

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1918.8572
$ws.Range("I81").Value = 942.5
$ws.Range("K81").Value = 1885
$ws.Range("M81").Value = -824
$ws.Range("H84").Value = 1918.8572
$ws.Range("I84").Value = 942.5
$ws.Range("K84").Value = 9425
$ws.Range("M84").Value = -4121
$ws.Range("H100").Value = 817.125
$ws.Range("J100").Value = 999
$ws.Range("L100").Value = 1998
$ws.Range("N100").Value = -3080
$ws.Range("H103").Value = 44101
$ws.Range("J103").Value = 44101
$ws.Range("L103").Value = 44101
$ws.Range("N103").Value = -46445
$ws.Range("H105").Value = 70000
$ws.Range("J105").Value = 70000
$ws.Range("L105").Value = 70000
$ws.Range("N105").Value = -76988
$ws.Range("H122").Value = 4033.875
$ws.Range("I122").Value = 2378.5
$ws.Range("K122").Value = 7135.5
$ws.Range("M122").Value = -4685.5
$ws.Range("H125").Value = 77578.5
$ws.Range("J125").Value = 77578.5
$ws.Range("L125").Value = 77578.5
$ws.Range("N125").Value = -87418.5
$ws.Range("H132").Value = 2898.6597
$ws.Range("I132").Value = 3184.6943
$ws.Range("K132").Value = 9554.082900000001
$ws.Range("M132").Value = -7024.082900000001
$ws.Range("H136").Value = 4800.4136
$ws.Range("I136").Value = 2296.0476
$ws.Range("K136").Value = 6888.1428
$ws.Range("M136").Value = -4338.1428
